$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row data corrections (rows were mismatched/out of order; swapping back to correct pairing) ---

# Row 10
$ws.Cells.Item(10, 2).Value2 = 6149376
$ws.Cells.Item(10, 6).Value2 = 'Jeonbuk Motors'
$ws.Cells.Item(10, 7).Value2 = 'Daegu FC'
$ws.Cells.Item(10, 8).Value2 = 1
$ws.Cells.Item(10, 9).Value2 = 0
$ws.Cells.Item(10, 10).Value2 = 'H'
$ws.Cells.Item(10, 11).Value2 = 1.75
$ws.Cells.Item(10, 12).Value2 = 3.5
$ws.Cells.Item(10, 13).Value2 = 4.2
$ws.Cells.Item(10, 14).Value2 = 2.05
$ws.Cells.Item(10, 15).Value2 = 3.3
$ws.Cells.Item(10, 16).Value2 = 3.4
$ws.Cells.Item(10, 17).Value2 = -0.25
$ws.Cells.Item(10, 18).Value2 = 1.825
$ws.Cells.Item(10, 19).Value2 = 2.025
$ws.Cells.Item(10, 20).Value2 = 2.25
$ws.Cells.Item(10, 21).Value2 = 1.85
$ws.Cells.Item(10, 22).Value2 = 2
$ws.Cells.Item(10, 23).Value2 = 1.05
$ws.Cells.Item(10, 24).Value2 = -1
$ws.Cells.Item(10, 25).Value2 = -1
$ws.Cells.Item(10, 26).Value2 = 0.825
$ws.Cells.Item(10, 27).Value2 = -1
$ws.Cells.Item(10, 28).Value2 = -1
$ws.Cells.Item(10, 29).Value2 = 1

# Row 11
$ws.Cells.Item(11, 2).Value2 = 6149854
$ws.Cells.Item(11, 6).Value2 = 'Gwangju FC'
$ws.Cells.Item(11, 7).Value2 = 'Suwon Bluewings'
$ws.Cells.Item(11, 8).Value2 = 2
$ws.Cells.Item(11, 9).Value2 = 1
$ws.Cells.Item(11, 10).Value2 = 'H'
$ws.Cells.Item(11, 11).Value2 = 1.8
$ws.Cells.Item(11, 12).Value2 = 3.4
$ws.Cells.Item(11, 13).Value2 = 3.8
$ws.Cells.Item(11, 14).Value2 = 1.85
$ws.Cells.Item(11, 15).Value2 = 3.6
$ws.Cells.Item(11, 16).Value2 = 4.2
$ws.Cells.Item(11, 17).Value2 = -0.5
$ws.Cells.Item(11, 18).Value2 = 1.9
$ws.Cells.Item(11, 19).Value2 = 1.95
$ws.Cells.Item(11, 20).Value2 = 2.5
$ws.Cells.Item(11, 21).Value2 = 2.025
$ws.Cells.Item(11, 22).Value2 = 1.825
$ws.Cells.Item(11, 23).Value2 = 0.8500000000000001
$ws.Cells.Item(11, 24).Value2 = -1
$ws.Cells.Item(11, 25).Value2 = -1
$ws.Cells.Item(11, 26).Value2 = 0.8999999999999999
$ws.Cells.Item(11, 27).Value2 = -1
$ws.Cells.Item(11, 28).Value2 = 1.025
$ws.Cells.Item(11, 29).Value2 = -1

# Row 41
$ws.Cells.Item(41, 2).Value2 = 6149871
$ws.Cells.Item(41, 6).Value2 = 'Suwon Bluewings'
$ws.Cells.Item(41, 7).Value2 = 'Pohang Steelers'
$ws.Cells.Item(41, 8).Value2 = 1
$ws.Cells.Item(41, 9).Value2 = 1
$ws.Cells.Item(41, 10).Value2 = 'D'
$ws.Cells.Item(41, 11).Value2 = 4
$ws.Cells.Item(41, 12).Value2 = 3.3
$ws.Cells.Item(41, 13).Value2 = 1.909
$ws.Cells.Item(41, 14).Value2 = 4.5
$ws.Cells.Item(41, 15).Value2 = 3.3
$ws.Cells.Item(41, 16).Value2 = 1.833
$ws.Cells.Item(41, 17).Value2 = 0.5
$ws.Cells.Item(41, 18).Value2 = 2.025
$ws.Cells.Item(41, 19).Value2 = 1.825
$ws.Cells.Item(41, 20).Value2 = 2.25
$ws.Cells.Item(41, 21).Value2 = 1.925
$ws.Cells.Item(41, 22).Value2 = 1.925
$ws.Cells.Item(41, 23).Value2 = -1
$ws.Cells.Item(41, 24).Value2 = 2.3
$ws.Cells.Item(41, 25).Value2 = -1
$ws.Cells.Item(41, 26).Value2 = 1.025
$ws.Cells.Item(41, 27).Value2 = -1
$ws.Cells.Item(41, 28).Value2 = -0.5
$ws.Cells.Item(41, 29).Value2 = 0.4625

# Row 42
$ws.Cells.Item(42, 2).Value2 = 6149872
$ws.Cells.Item(42, 6).Value2 = 'FC Seoul'
$ws.Cells.Item(42, 7).Value2 = 'Suwon FC'
$ws.Cells.Item(42, 8).Value2 = 7
$ws.Cells.Item(42, 9).Value2 = 2
$ws.Cells.Item(42, 10).Value2 = 'H'
$ws.Cells.Item(42, 11).Value2 = 1.55
$ws.Cells.Item(42, 12).Value2 = 3.75
$ws.Cells.Item(42, 13).Value2 = 5.25
$ws.Cells.Item(42, 14).Value2 = 1.55
$ws.Cells.Item(42, 15).Value2 = 4
$ws.Cells.Item(42, 16).Value2 = 5
$ws.Cells.Item(42, 17).Value2 = -1
$ws.Cells.Item(42, 18).Value2 = 1.925
$ws.Cells.Item(42, 19).Value2 = 1.925
$ws.Cells.Item(42, 20).Value2 = 3
$ws.Cells.Item(42, 21).Value2 = 1.9
$ws.Cells.Item(42, 22).Value2 = 1.95
$ws.Cells.Item(42, 23).Value2 = 0.55
$ws.Cells.Item(42, 24).Value2 = -1
$ws.Cells.Item(42, 25).Value2 = -1
$ws.Cells.Item(42, 26).Value2 = 0.925
$ws.Cells.Item(42, 27).Value2 = -1
$ws.Cells.Item(42, 28).Value2 = 0.8999999999999999
$ws.Cells.Item(42, 29).Value2 = -1

# Row 43
$ws.Cells.Item(43, 2).Value2 = 6149381
$ws.Cells.Item(43, 6).Value2 = 'Daejeon Hana Citizen'
$ws.Cells.Item(43, 7).Value2 = 'Jeonbuk Motors'
$ws.Cells.Item(43, 8).Value2 = 2
$ws.Cells.Item(43, 9).Value2 = 2
$ws.Cells.Item(43, 10).Value2 = 'D'
$ws.Cells.Item(43, 11).Value2 = 3.4
$ws.Cells.Item(43, 12).Value2 = 3.25
$ws.Cells.Item(43, 13).Value2 = 2.05
$ws.Cells.Item(43, 14).Value2 = 3.25
$ws.Cells.Item(43, 15).Value2 = 3.3
$ws.Cells.Item(43, 16).Value2 = 2.1
$ws.Cells.Item(43, 17).Value2 = 0.25
$ws.Cells.Item(43, 18).Value2 = 1.975
$ws.Cells.Item(43, 19).Value2 = 1.875
$ws.Cells.Item(43, 20).Value2 = 2.75
$ws.Cells.Item(43, 21).Value2 = 2
$ws.Cells.Item(43, 22).Value2 = 1.85
$ws.Cells.Item(43, 23).Value2 = -1
$ws.Cells.Item(43, 24).Value2 = 2.3
$ws.Cells.Item(43, 25).Value2 = -1
$ws.Cells.Item(43, 26).Value2 = 0.4875
$ws.Cells.Item(43, 27).Value2 = -0.5
$ws.Cells.Item(43, 28).Value2 = 1
$ws.Cells.Item(43, 29).Value2 = -1

# Row 56
$ws.Cells.Item(56, 2).Value2 = 6149880
$ws.Cells.Item(56, 6).Value2 = 'FC Seoul'
$ws.Cells.Item(56, 7).Value2 = 'Pohang Steelers'
$ws.Cells.Item(56, 8).Value2 = 2
$ws.Cells.Item(56, 9).Value2 = 2
$ws.Cells.Item(56, 10).Value2 = 'D'
$ws.Cells.Item(56, 11).Value2 = 2.5
$ws.Cells.Item(56, 12).Value2 = 3.2
$ws.Cells.Item(56, 13).Value2 = 2.8
$ws.Cells.Item(56, 14).Value2 = 2.625
$ws.Cells.Item(56, 15).Value2 = 3.3
$ws.Cells.Item(56, 16).Value2 = 2.6
$ws.Cells.Item(56, 17).Value2 = 0
$ws.Cells.Item(56, 18).Value2 = 1.95
$ws.Cells.Item(56, 19).Value2 = 1.9
$ws.Cells.Item(56, 20).Value2 = 2.25
$ws.Cells.Item(56, 21).Value2 = 1.825
$ws.Cells.Item(56, 22).Value2 = 2.025
$ws.Cells.Item(56, 23).Value2 = -1
$ws.Cells.Item(56, 24).Value2 = 2.3
$ws.Cells.Item(56, 25).Value2 = -1
$ws.Cells.Item(56, 26).Value2 = 0
$ws.Cells.Item(56, 27).Value2 = 0
$ws.Cells.Item(56, 28).Value2 = 0.825
$ws.Cells.Item(56, 29).Value2 = -1

# Row 57
$ws.Cells.Item(57, 2).Value2 = 6149879
$ws.Cells.Item(57, 6).Value2 = 'Gwangju FC'
$ws.Cells.Item(57, 7).Value2 = 'Daejeon Hana Citizen'
$ws.Cells.Item(57, 8).Value2 = 3
$ws.Cells.Item(57, 9).Value2 = 0
$ws.Cells.Item(57, 10).Value2 = 'H'
$ws.Cells.Item(57, 11).Value2 = 2.15
$ws.Cells.Item(57, 12).Value2 = 3.3
$ws.Cells.Item(57, 13).Value2 = 3.4
$ws.Cells.Item(57, 14).Value2 = 2.2
$ws.Cells.Item(57, 15).Value2 = 3
$ws.Cells.Item(57, 16).Value2 = 3.6
$ws.Cells.Item(57, 17).Value2 = -0.25
$ws.Cells.Item(57, 18).Value2 = 1.9
$ws.Cells.Item(57, 19).Value2 = 1.95
$ws.Cells.Item(57, 20).Value2 = 2.25
$ws.Cells.Item(57, 21).Value2 = 2
$ws.Cells.Item(57, 22).Value2 = 1.85
$ws.Cells.Item(57, 23).Value2 = 1.2
$ws.Cells.Item(57, 24).Value2 = -1
$ws.Cells.Item(57, 25).Value2 = -1
$ws.Cells.Item(57, 26).Value2 = 0.8999999999999999
$ws.Cells.Item(57, 27).Value2 = -1
$ws.Cells.Item(57, 28).Value2 = 1
$ws.Cells.Item(57, 29).Value2 = -1

# Row 94
$ws.Cells.Item(94, 2).Value2 = 6323586
$ws.Cells.Item(94, 6).Value2 = 'Suwon FC'
$ws.Cells.Item(94, 7).Value2 = 'Ulsan Hyundai'
$ws.Cells.Item(94, 8).Value2 = 2
$ws.Cells.Item(94, 9).Value2 = 3
$ws.Cells.Item(94, 10).Value2 = 'A'
$ws.Cells.Item(94, 11).Value2 = 4.5
$ws.Cells.Item(94, 12).Value2 = 4
$ws.Cells.Item(94, 13).Value2 = 1.7
$ws.Cells.Item(94, 14).Value2 = 4.75
$ws.Cells.Item(94, 15).Value2 = 4
$ws.Cells.Item(94, 16).Value2 = 1.666
$ws.Cells.Item(94, 17).Value2 = 0.75
$ws.Cells.Item(94, 18).Value2 = 2.025
$ws.Cells.Item(94, 19).Value2 = 1.825
$ws.Cells.Item(94, 20).Value2 = 3.25
$ws.Cells.Item(94, 21).Value2 = 2.05
$ws.Cells.Item(94, 22).Value2 = 1.8
$ws.Cells.Item(94, 23).Value2 = -1
$ws.Cells.Item(94, 24).Value2 = -1
$ws.Cells.Item(94, 25).Value2 = 0.6659999999999999
$ws.Cells.Item(94, 26).Value2 = -0.5
$ws.Cells.Item(94, 27).Value2 = 0.4125
$ws.Cells.Item(94, 28).Value2 = 1.05
$ws.Cells.Item(94, 29).Value2 = -1

# Row 95
$ws.Cells.Item(95, 2).Value2 = 6323587
$ws.Cells.Item(95, 6).Value2 = 'Gwangju FC'
$ws.Cells.Item(95, 7).Value2 = 'Jeonbuk Motors'
$ws.Cells.Item(95, 8).Value2 = 0
$ws.Cells.Item(95, 9).Value2 = 1
$ws.Cells.Item(95, 10).Value2 = 'A'
$ws.Cells.Item(95, 11).Value2 = 2.7
$ws.Cells.Item(95, 12).Value2 = 3.4
$ws.Cells.Item(95, 13).Value2 = 2.45
$ws.Cells.Item(95, 14).Value2 = 2.6
$ws.Cells.Item(95, 15).Value2 = 3.3
$ws.Cells.Item(95, 16).Value2 = 2.7
$ws.Cells.Item(95, 17).Value2 = 0
$ws.Cells.Item(95, 18).Value2 = 1.875
$ws.Cells.Item(95, 19).Value2 = 1.975
$ws.Cells.Item(95, 20).Value2 = 2.25
$ws.Cells.Item(95, 21).Value2 = 1.975
$ws.Cells.Item(95, 22).Value2 = 1.875
$ws.Cells.Item(95, 23).Value2 = -1
$ws.Cells.Item(95, 24).Value2 = -1
$ws.Cells.Item(95, 25).Value2 = 1.7
$ws.Cells.Item(95, 26).Value2 = -1
$ws.Cells.Item(95, 27).Value2 = 0.9750000000000001
$ws.Cells.Item(95, 28).Value2 = -1
$ws.Cells.Item(95, 29).Value2 = 0.875

# Row 129
$ws.Cells.Item(129, 2).Value2 = 7333495
$ws.Cells.Item(129, 6).Value2 = 'Jeonbuk Motors'
$ws.Cells.Item(129, 7).Value2 = 'Gwangju FC'
$ws.Cells.Item(129, 8).Value2 = 2
$ws.Cells.Item(129, 9).Value2 = 0
$ws.Cells.Item(129, 10).Value2 = 'H'
$ws.Cells.Item(129, 11).Value2 = 2.05
$ws.Cells.Item(129, 12).Value2 = 3.25
$ws.Cells.Item(129, 13).Value2 = 3.2
$ws.Cells.Item(129, 14).Value2 = 2.05
$ws.Cells.Item(129, 15).Value2 = 3.4
$ws.Cells.Item(129, 16).Value2 = 3.6
$ws.Cells.Item(129, 17).Value2 = -0.25
$ws.Cells.Item(129, 18).Value2 = 1.775
$ws.Cells.Item(129, 19).Value2 = 2.1
$ws.Cells.Item(129, 20).Value2 = 2.25
$ws.Cells.Item(129, 21).Value2 = 1.9
$ws.Cells.Item(129, 22).Value2 = 1.95
$ws.Cells.Item(129, 23).Value2 = 1.05
$ws.Cells.Item(129, 24).Value2 = -1
$ws.Cells.Item(129, 25).Value2 = -1
$ws.Cells.Item(129, 26).Value2 = 0.7749999999999999
$ws.Cells.Item(129, 27).Value2 = -1
$ws.Cells.Item(129, 28).Value2 = -0.5
$ws.Cells.Item(129, 29).Value2 = 0.475

# Row 130
$ws.Cells.Item(130, 2).Value2 = 7334087
$ws.Cells.Item(130, 6).Value2 = 'Jeju United'
$ws.Cells.Item(130, 7).Value2 = 'Daejeon Hana Citizen'
$ws.Cells.Item(130, 8).Value2 = 0
$ws.Cells.Item(130, 9).Value2 = 2
$ws.Cells.Item(130, 10).Value2 = 'A'
$ws.Cells.Item(130, 11).Value2 = 2.05
$ws.Cells.Item(130, 12).Value2 = 3.75
$ws.Cells.Item(130, 13).Value2 = 3.25
$ws.Cells.Item(130, 14).Value2 = 2.05
$ws.Cells.Item(130, 15).Value2 = 3.75
$ws.Cells.Item(130, 16).Value2 = 3.3
$ws.Cells.Item(130, 17).Value2 = -0.25
$ws.Cells.Item(130, 18).Value2 = 1.8
$ws.Cells.Item(130, 19).Value2 = 2.05
$ws.Cells.Item(130, 20).Value2 = 2.75
$ws.Cells.Item(130, 21).Value2 = 1.975
$ws.Cells.Item(130, 22).Value2 = 1.875
$ws.Cells.Item(130, 23).Value2 = -1
$ws.Cells.Item(130, 24).Value2 = -1
$ws.Cells.Item(130, 25).Value2 = 2.3
$ws.Cells.Item(130, 26).Value2 = -1
$ws.Cells.Item(130, 27).Value2 = 1.05
$ws.Cells.Item(130, 28).Value2 = -1
$ws.Cells.Item(130, 29).Value2 = 0.875

# Row 138
$ws.Cells.Item(138, 2).Value2 = 7333499
$ws.Cells.Item(138, 6).Value2 = 'Daegu FC'
$ws.Cells.Item(138, 7).Value2 = 'Incheon Utd'
$ws.Cells.Item(138, 8).Value2 = 2
$ws.Cells.Item(138, 9).Value2 = 1
$ws.Cells.Item(138, 10).Value2 = 'H'
$ws.Cells.Item(138, 11).Value2 = 2.55
$ws.Cells.Item(138, 12).Value2 = 3.3
$ws.Cells.Item(138, 13).Value2 = 2.7
$ws.Cells.Item(138, 14).Value2 = 2.8
$ws.Cells.Item(138, 15).Value2 = 3.3
$ws.Cells.Item(138, 16).Value2 = 2.45
$ws.Cells.Item(138, 17).Value2 = 0
$ws.Cells.Item(138, 18).Value2 = 2.1
$ws.Cells.Item(138, 19).Value2 = 1.775
$ws.Cells.Item(138, 20).Value2 = 2.25
$ws.Cells.Item(138, 21).Value2 = 1.875
$ws.Cells.Item(138, 22).Value2 = 1.975
$ws.Cells.Item(138, 23).Value2 = 1.8
$ws.Cells.Item(138, 24).Value2 = -1
$ws.Cells.Item(138, 25).Value2 = -1
$ws.Cells.Item(138, 26).Value2 = 1.1
$ws.Cells.Item(138, 27).Value2 = -1
$ws.Cells.Item(138, 28).Value2 = 0.875
$ws.Cells.Item(138, 29).Value2 = -1

# Row 139
$ws.Cells.Item(139, 2).Value2 = 7333498
$ws.Cells.Item(139, 6).Value2 = 'Gwangju FC'
$ws.Cells.Item(139, 7).Value2 = 'Pohang Steelers'
$ws.Cells.Item(139, 8).Value2 = 0
$ws.Cells.Item(139, 9).Value2 = 0
$ws.Cells.Item(139, 10).Value2 = 'D'
$ws.Cells.Item(139, 11).Value2 = 2
$ws.Cells.Item(139, 12).Value2 = 3.4
$ws.Cells.Item(139, 13).Value2 = 3.3
$ws.Cells.Item(139, 14).Value2 = 1.8
$ws.Cells.Item(139, 15).Value2 = 3.5
$ws.Cells.Item(139, 16).Value2 = 4
$ws.Cells.Item(139, 17).Value2 = -0.5
$ws.Cells.Item(139, 18).Value2 = 1.8
$ws.Cells.Item(139, 19).Value2 = 2.05
$ws.Cells.Item(139, 20).Value2 = 2.5
$ws.Cells.Item(139, 21).Value2 = 2.025
$ws.Cells.Item(139, 22).Value2 = 1.825
$ws.Cells.Item(139, 23).Value2 = -1
$ws.Cells.Item(139, 24).Value2 = 2.5
$ws.Cells.Item(139, 25).Value2 = -1
$ws.Cells.Item(139, 26).Value2 = -1
$ws.Cells.Item(139, 27).Value2 = 1.05
$ws.Cells.Item(139, 28).Value2 = -1
$ws.Cells.Item(139, 29).Value2 = 0.825

# Row 143
$ws.Cells.Item(143, 2).Value2 = 7715259
$ws.Cells.Item(143, 6).Value2 = 'Incheon Utd'
$ws.Cells.Item(143, 7).Value2 = 'Suwon FC'
$ws.Cells.Item(143, 8).Value2 = 0
$ws.Cells.Item(143, 9).Value2 = 1
$ws.Cells.Item(143, 10).Value2 = 'A'
$ws.Cells.Item(143, 11).Value2 = 2
$ws.Cells.Item(143, 12).Value2 = 3.5
$ws.Cells.Item(143, 13).Value2 = 3.5
$ws.Cells.Item(143, 14).Value2 = 1.909
$ws.Cells.Item(143, 15).Value2 = 3.4
$ws.Cells.Item(143, 16).Value2 = 3.75
$ws.Cells.Item(143, 17).Value2 = -0.25
$ws.Cells.Item(143, 18).Value2 = 1.8
$ws.Cells.Item(143, 19).Value2 = 2.05
$ws.Cells.Item(143, 20).Value2 = 2.5
$ws.Cells.Item(143, 21).Value2 = 1.95
$ws.Cells.Item(143, 22).Value2 = 1.9
$ws.Cells.Item(143, 23).Value2 = -1
$ws.Cells.Item(143, 24).Value2 = -1
$ws.Cells.Item(143, 25).Value2 = 2.75
$ws.Cells.Item(143, 26).Value2 = -1
$ws.Cells.Item(143, 27).Value2 = 1.05
$ws.Cells.Item(143, 28).Value2 = -1
$ws.Cells.Item(143, 29).Value2 = 0.8999999999999999

# Row 144
$ws.Cells.Item(144, 2).Value2 = 7716460
$ws.Cells.Item(144, 6).Value2 = 'Gangwon FC'
$ws.Cells.Item(144, 7).Value2 = 'Jeju United'
$ws.Cells.Item(144, 8).Value2 = 1
$ws.Cells.Item(144, 9).Value2 = 1
$ws.Cells.Item(144, 10).Value2 = 'D'
$ws.Cells.Item(144, 11).Value2 = 2.5
$ws.Cells.Item(144, 12).Value2 = 3
$ws.Cells.Item(144, 13).Value2 = 3
$ws.Cells.Item(144, 14).Value2 = 2.875
$ws.Cells.Item(144, 15).Value2 = 2.7
$ws.Cells.Item(144, 16).Value2 = 2.9
$ws.Cells.Item(144, 17).Value2 = 0
$ws.Cells.Item(144, 18).Value2 = 1.925
$ws.Cells.Item(144, 19).Value2 = 1.925
$ws.Cells.Item(144, 20).Value2 = 1.75
$ws.Cells.Item(144, 21).Value2 = 1.9
$ws.Cells.Item(144, 22).Value2 = 1.95
$ws.Cells.Item(144, 23).Value2 = -1
$ws.Cells.Item(144, 24).Value2 = 1.7
$ws.Cells.Item(144, 25).Value2 = -1
$ws.Cells.Item(144, 26).Value2 = 0
$ws.Cells.Item(144, 27).Value2 = 0
$ws.Cells.Item(144, 28).Value2 = 0.45
$ws.Cells.Item(144, 29).Value2 = -0.5

# Row 147
$ws.Cells.Item(147, 2).Value2 = 7715262
$ws.Cells.Item(147, 6).Value2 = 'Gimcheon Sangmu FC'
$ws.Cells.Item(147, 7).Value2 = 'Ulsan Hyundai'
$ws.Cells.Item(147, 8).Value2 = 2
$ws.Cells.Item(147, 9).Value2 = 3
$ws.Cells.Item(147, 10).Value2 = 'A'
$ws.Cells.Item(147, 11).Value2 = 3.3
$ws.Cells.Item(147, 12).Value2 = 3.5
$ws.Cells.Item(147, 13).Value2 = 2.05
$ws.Cells.Item(147, 14).Value2 = 3
$ws.Cells.Item(147, 15).Value2 = 3.2
$ws.Cells.Item(147, 16).Value2 = 2.4
$ws.Cells.Item(147, 17).Value2 = 0.25
$ws.Cells.Item(147, 18).Value2 = 1.775
$ws.Cells.Item(147, 19).Value2 = 2.1
$ws.Cells.Item(147, 20).Value2 = 2.25
$ws.Cells.Item(147, 21).Value2 = 2
$ws.Cells.Item(147, 22).Value2 = 1.85
$ws.Cells.Item(147, 23).Value2 = -1
$ws.Cells.Item(147, 24).Value2 = -1
$ws.Cells.Item(147, 25).Value2 = 1.4
$ws.Cells.Item(147, 26).Value2 = -1
$ws.Cells.Item(147, 27).Value2 = 1.1
$ws.Cells.Item(147, 28).Value2 = 1
$ws.Cells.Item(147, 29).Value2 = -1

# Row 148
$ws.Cells.Item(148, 2).Value2 = 7715261
$ws.Cells.Item(148, 6).Value2 = 'Pohang Steelers'
$ws.Cells.Item(148, 7).Value2 = 'Daegu FC'
$ws.Cells.Item(148, 8).Value2 = 3
$ws.Cells.Item(148, 9).Value2 = 1
$ws.Cells.Item(148, 10).Value2 = 'H'
$ws.Cells.Item(148, 11).Value2 = 2.05
$ws.Cells.Item(148, 12).Value2 = 3.3
$ws.Cells.Item(148, 13).Value2 = 3.6
$ws.Cells.Item(148, 14).Value2 = 2.75
$ws.Cells.Item(148, 15).Value2 = 3.1
$ws.Cells.Item(148, 16).Value2 = 2.7
$ws.Cells.Item(148, 17).Value2 = 0
$ws.Cells.Item(148, 18).Value2 = 1.95
$ws.Cells.Item(148, 19).Value2 = 1.9
$ws.Cells.Item(148, 20).Value2 = 2
$ws.Cells.Item(148, 21).Value2 = 1.85
$ws.Cells.Item(148, 22).Value2 = 2
$ws.Cells.Item(148, 23).Value2 = 1.75
$ws.Cells.Item(148, 24).Value2 = -1
$ws.Cells.Item(148, 25).Value2 = -1
$ws.Cells.Item(148, 26).Value2 = 0.95
$ws.Cells.Item(148, 27).Value2 = -1
$ws.Cells.Item(148, 28).Value2 = 0.8500000000000001
$ws.Cells.Item(148, 29).Value2 = -1

# --- Direct odds corrections for upcoming fixtures (rows 158-160) ---

# Row 158
$ws.Cells.Item(158, 14).Value2 = 2
$ws.Cells.Item(158, 16).Value2 = 3.4
$ws.Cells.Item(158, 17).Value2 = -0.25
$ws.Cells.Item(158, 18).Value2 = 1.825
$ws.Cells.Item(158, 19).Value2 = 2.025

# Row 159
$ws.Cells.Item(159, 18).Value2 = 1.8
$ws.Cells.Item(159, 19).Value2 = 2.05

# Row 160
$ws.Cells.Item(160, 18).Value2 = 1.825
$ws.Cells.Item(160, 19).Value2 = 2.025
